{"js": "// Added description on report:\n//   - Extend / rewrite the tail of the \"dataset\" paragraph with a new\n//     sentence describing the stakeholders the analysis targets.\n//   - Remove the now-superfluous empty trailing paragraph.\n\nconst body = context.document.body;\n\n// 1) Replace the old closing sentence of the \"dataset\" paragraph with the\n//    new description of who the analysis targets.\nconst oldTail =\n  \"helps both the public office and private business pointing out where \" +\n  \"the infrastructure is lacking so that, a newer infrastructure with a \" +\n  \"nominal service fee would be able to generate a large revenue.\";\nconst newTail =\n  \"targets multiple stakeholders like Public offices of the Los Angeles, \" +\n  \"The LAPD and also private business who might benefit from making \" +\n  \"additional private pay-to-park complex in the vicinity.\";\n\nconst found = body.search(oldTail, { matchCase: true });\nfound.load(\"items\");\nawait context.sync();\n\nif (found.items.length > 0) {\n  // Replace in place so the surrounding run formatting (bold/size, etc.)\n  // is preserved.\n  found.items[0].insertText(newTail, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Drop the trailing empty paragraph left at the end of the document.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst count = paragraphs.items.length;\nif (count > 1) {\n  const lastParagraph = paragraphs.items[count - 1];\n  lastParagraph.load(\"text\");\n  await context.sync();\n\n  if (lastParagraph.text.trim().length === 0) {\n    // Paragraph.delete() on the very last paragraph of the body is a\n    // no-op (Word always needs a final paragraph mark), so instead we\n    // remove the range spanning from the end of the previous paragraph\n    // through the end of this one, which deletes its paragraph mark too.\n    const previousParagraph = paragraphs.items[count - 2];\n    const rangeStart = previousParagraph.getRange(\"End\");\n    const rangeEnd = lastParagraph.getRange(\"End\");\n    const rangeToDelete = rangeStart.expandTo(rangeEnd);\n    rangeToDelete.delete();\n    await context.sync();\n  }\n}\n", "ps1": "# Added description on report:\n#   - Extend / rewrite the tail of the \"dataset\" paragraph with a new\n#     sentence describing the stakeholders the analysis targets.\n#   - Remove the now-superfluous empty trailing paragraph.\n\n$d = $word.ActiveDocument\n\n# 1) Replace the old closing sentence of the \"dataset\" paragraph with the\n#    new description of who the analysis targets. Using Find/Replace keeps\n#    the surrounding run formatting (bold/size, etc.) intact.\n$oldTail = \"helps both the public office and private business pointing out where the infrastructure is lacking so that, a newer infrastructure with a nominal service fee would be able to generate a large revenue.\"\n$newTail = \"targets multiple stakeholders like Public offices of the Los Angeles, The LAPD and also private business who might benefit from making additional private pay-to-park complex in the vicinity.\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\n    $oldTail,   # FindText\n    $false,     # MatchCase\n    $false,     # MatchWholeWord\n    $false,     # MatchWildcards\n    $false,     # MatchSoundsLike\n    $false,     # MatchAllWordForms\n    $true,      # Forward\n    1,          # Wrap (wdFindContinue)\n    $false,     # Format\n    $newTail,   # ReplaceWith\n    2           # Replace (wdReplaceAll)\n) | Out-Null\n\n# 2) Drop the trailing empty paragraph left at the end of the document.\n$count = $d.Paragraphs.Count\nif ($count -gt 1) {\n    $lastParagraph = $d.Paragraphs.Item($count)\n    if ($lastParagraph.Range.Text.Trim().Length -eq 0) {\n        # The very last paragraph mark of the document body can't be\n        # deleted directly, so instead remove the paragraph mark that\n        # separates the previous paragraph from this (now empty) one;\n        # that merges them away, leaving this paragraph's own mark as\n        # the new (required) final mark of the document.\n        $markStart = $lastParagraph.Range.Start - 1\n        $markEnd = $lastParagraph.Range.Start\n        $d.Range($markStart, $markEnd).Delete()\n    }\n}\n"}
